$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.843.21'
$ws.Range("E2").Value = '  +0.60%  '

# Row 3
$ws.Range("D3").Value = '1.641.67'
$ws.Range("E3").Value = '  +0.61%  '

# Row 4
$ws.Range("E4").Value = '  -0.79%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.99%  '

# Row 6
$ws.Range("E6").Value = '  +1.74%  '

# Row 7
$ws.Range("E7").Value = '  -0.75%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.253'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.86%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0621'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.36%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.75'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.14%  '

# Row 11
$ws.Range("E11").Value = '  +0.25%  '

# Row 12
$ws.Range("D12").Value = '1.870.88'
$ws.Range("E12").Value = '  +0.59%  '

# Row 13
$ws.Range("D13").Value = '1.631.02'
$ws.Range("E13").Value = '  -0.17%  '

# Row 14
$ws.Range("E14").Value = '  +0.60%  '

# Row 15
$ws.Range("E15").Value = '  +1.17%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.37'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.98%  '

# Row 17
$ws.Range("D17").Value = '26.846.86'
$ws.Range("E17").Value = '  +0.74%  '

# Row 18
$ws.Range("E18").Value = '  +1.51%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '217.63'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.84%  '

# Row 20
$ws.Range("E20").Value = '  -0.70%  '

# Row 21
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.62'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.75%  '

# Row 22
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.37'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.58%  '

# Row 23
$ws.Range("E23").Value = '  +4.28%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.12%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.49'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.10%  '

# Row 26
$ws.Range("E26").Value = '  -0.98%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.37'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.26%  '

# Row 28
$ws.Range("E28").Value = '  +1.20%  '

# Row 29
$ws.Range("E29").Value = '  +2.27%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0511'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.31%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.36'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.69%  '

# Row 33
$ws.Range("E33").Value = '  +2.01%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.55'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.91%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.44'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.18%  '

# Row 36
$ws.Range("D36").Value = '1.237.57'
$ws.Range("E36").Value = '  -1.69%  '

# Row 37
$ws.Range("E37").Value = '  +0.60%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.539'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.75%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.833'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.58%  '

# Row 40
$ws.Range("E40").Value = '  -0.70%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.806'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.92%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.37'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.54%  '

# Row 43
$ws.Range("D43").Value = '1.783.14'
$ws.Range("E43").Value = '  +0.57%  '

# Row 44
$ws.Range("E44").Value = '  -2.51%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '60.71'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.04%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.46'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.72%  '

# Row 47
$ws.Range("E47").Value = '  +0.93%  '

# Row 48
$ws.Range("E48").Value = '  +4.20%  '

# Row 49
$ws.Range("E49").Value = '  -0.61%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0973'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.35%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.56'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.54%  '
